$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 51: Bruno Díaz / 2017-05-26 / 4 hours / Sprint 3 - Integración BackEnd y FrontEnd /
# Investigación sobre API REST - Servicios del back ya responden al llamarlos desde el front
$ws.Cells.Item(51, 1).Value = "Bruno Díaz"

# Set the date as a raw serial number, then copy the number formatting from the
# cell above (B50) so it reuses the existing date style instead of generating
# a brand-new (duplicate) number format.
$ws.Cells.Item(51, 2).Value = 42881
$ws.Cells.Item(50, 2).Copy()
$ws.Cells.Item(51, 2).PasteSpecial(-4122)

$ws.Cells.Item(51, 3).Value = 4
$ws.Cells.Item(51, 4).Value = "Sprint 3 - Integración BackEnd y FrontEnd"
$ws.Cells.Item(51, 5).Value = "Investigación sobre API REST - Servicios del back ya responden al llamarlos desde el front"

# Clear the marching-ants from the copy operation and move the active
# selection to E52, matching the saved cursor position in the workbook.
$excel.CutCopyMode = 0
$ws.Range("E52").Select() | Out-Null
